$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.133.85"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.902.50"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'252.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "'0.697"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'41.65"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").Value = "'0.353"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").Value = "'52.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "'0.0756"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("D12").Value = "'0.0978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'13.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.35%  "
$ws.Range("D14").Value = "2.180.26"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").Value = "'5.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("D17").Value = "1.913.95"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "35.151.12"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "'73.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").Value = "'242.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "'13.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "'5.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "'168.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "'8.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'18.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "'0.129"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "4.128.73"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'2.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.66%  "
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("D34").Value = "'0.0595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("D35").Value = "'1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.30%  "
$ws.Range("D36").Value = "'4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "'0.849"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.55%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'17.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("D41").Value = "'97.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("D42").Value = "'0.0215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("D43").Value = "'0.0663"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").Value = "'1.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "'2.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "1.304.20"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "'6.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D50").Value = "'11.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "'0.0753"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.15%  "
